$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F5").Value = 344
$ws.Range("F6").Value = 7825
$ws.Range("F9").Value = 1166
$ws.Range("F14").Value = 1797
$ws.Range("F22").Value = 8413
$ws.Range("F28").Value = 1051
$ws.Range("F29").Value = 543
$ws.Range("F30").Value = 1114
$ws.Range("F31").Value = 549
$ws.Range("F32").Value = 549
$ws.Range("F37").Value = 1016
$ws.Range("F38").Value = 611
$ws.Range("F41").Value = 3370
$ws.Range("F42").Value = 915
$ws.Range("F43").Value = 25
$ws.Range("F46").Value = 64
$ws.Range("F48").Value = 426
$ws.Range("F49").Value = 538

$ws = $wb.Worksheets.Item(2)
$ws.Range("G4").Value = 108
$ws.Range("F7").Value = 317
$ws.Range("F18").Value = 59
$ws.Range("F22").Value = 26
$ws.Range("F24").Value = 75
$ws.Range("F25").Value = 6787

$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 2031
$ws.Range("F8").Value = 2193
$ws.Range("F9").Value = 9014
$ws.Range("F10").Value = 1209

$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 2031
$ws.Range("F4").Value = 7825
$ws.Range("F7").Value = 2193
$ws.Range("F9").Value = 1209
$ws.Range("F11").Value = 1166
$ws.Range("F15").Value = 317
$ws.Range("F22").Value = 8413
$ws.Range("F27").Value = 1051
$ws.Range("F28").Value = 543
$ws.Range("F29").Value = 1114
$ws.Range("F30").Value = 549
$ws.Range("F31").Value = 549
$ws.Range("F34").Value = 1016
$ws.Range("F35").Value = 611
$ws.Range("F38").Value = 3370
$ws.Range("F39").Value = 915
$ws.Range("F43").Value = 426
$ws.Range("F44").Value = 538
